# Apply updated crypto price/volume data per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (the sheet stores these as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Write the new values.
$ws.Range("D2").Value = "51.963.73"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "3.029.70"
$ws.Range("E3").Value = "  +4.23%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "382.85"
$ws.Range("E5").Value = "  +5.92%  "
$ws.Range("D6").Value = "107.07"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D7").Value = "0.551"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").Value = "38.19"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "0.0852"
$ws.Range("E12").Value = "  +2.22%  "
$ws.Range("D13").Value = "18.99"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("D14").Value = "3.494.03"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "7.60"
$ws.Range("E15").Value = "  +3.58%  "
$ws.Range("D16").Value = "3.010.22"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").Value = "0.984"
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("D18").Value = "51.967.68"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "3.45"
$ws.Range("E19").Value = "  +4.99%  "
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  +4.28%  "
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").Value = "0.0₃0970"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "69.21"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").Value = "265.46"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("E25").Value = "  +4.63%  "
$ws.Range("D26").Value = "0.174"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "7.34"
$ws.Range("E27").Value = "  +19.39%  "
$ws.Range("D28").Value = "7.60"
$ws.Range("E28").Value = "  +5.23%  "
$ws.Range("D29").Value = "26.42"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "35.36"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("E35").Value = "  -2.60%  "
$ws.Range("D36").Value = "0.0452"
$ws.Range("E36").Value = "  +7.29%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  +0.35%  "
$ws.Range("D39").Value = "17.77"
$ws.Range("E39").Value = "  +4.71%  "
$ws.Range("D40").Value = "2.68"
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +3.69%  "
$ws.Range("D43").Value = "125.27"
$ws.Range("E43").Value = "  +5.17%  "
$ws.Range("D44").Value = "22.66"
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "0.279"
$ws.Range("E46").Value = "  +18.55%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("D48").Value = "2.068.66"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "3.32"
$ws.Range("E49").Value = "  +4.32%  "
$ws.Range("D50").Value = "0.0356"
$ws.Range("E50").Value = "  +15.68%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "5.27"
$ws.Range("E51").Value = "  +4.69%  "
